# Apply "custom accuracy" rounding to row 5 values and remove row 6
# (dataset trimmed from 6 rows to 5 rows of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: round the measurement columns (B:AH) to the new custom accuracy.
$row5 = @{
    "B5"  = 16.21
    "C5"  = 12.11
    "D5"  = 1.06
    "E5"  = 35.54
    "F5"  = 28.87
    "G5"  = 12.43
    "H5"  = 51.55
    "I5"  = 19.71
    "J5"  = 8.93
    "K5"  = 12.74
    "L5"  = 14.24
    "M5"  = 15.2
    "N5"  = 4.27
    "O5"  = 12.78
    "P5"  = 18.08
    "Q5"  = 10.89
    "R5"  = 0.47
    "S5"  = 0.63
    "T5"  = 187.25
    "U5"  = 35.8
    "V5"  = 11.8
    "W5"  = 23.95
    "X5"  = 12.61
    "Y5"  = 1.66
    "Z5"  = 25.36
    "AA5" = 10.42
    "AB5" = 9.279999999999999
    "AC5" = 10.88
    "AD5" = 15.02
    "AE5" = 0.48
    "AF5" = 46.99
    "AG5" = 6.58
    "AH5" = 14.75
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove row 6 entirely (the last data row of the previous revision).
$ws.Rows.Item(6).Delete()

Write-Output "applied custom-accuracy rounding and dropped row 6"
